# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1. Adds a new "ODI Bowling Extra" worksheet (after "ODI Batting Extra")
#    containing MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns.
# 2. Cleans up stray empty placeholder cells (B2:E2) left on the
#    "ODI Batting Extra" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Remove the stray empty cells B2:E2 on "ODI Batting Extra" (row for
#    match 4471 never had batting-position/4s/6s/runs% data scraped).
# ---------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$battingExtra.Range("B2:E2").ClearContents()

# ---------------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" sheet at the end of the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row
$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Match the bold/bordered header look used on the other sheets.
$battingExtra.Range("A1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

# Data rows (MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL)
$data = @(
  @("4469", "0", "10.00%"),
  @("4470", "0", "10.00%"),
  @("4480", "0", ""),
  @("4482", "1", ""),
  @("4485", "0", "10.00%"),
  @("4487", "0", "10.00%"),
  @("4488", "0", "20.00%"),
  @("4521", "1", "30.00%"),
  @("4523", "", ""),
  @("4527", "0", "10.00%"),
  @("4594", "", ""),
  @("4597", "0", "30.00%"),
  @("4600", "", ""),
  @("4601", "1", "20.00%"),
  @("4603", "0", ""),
  @("4687", "", ""),
  @("4689", "0", "20.00%"),
  @("4691", "", ""),
  @("4735", "0", "40.00%"),
  @("4745", "", "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $i + 2

    $code = $rowValues[0]
    $maidenOvers = $rowValues[1]
    $percentWickets = $rowValues[2]

    $bowlingExtra.Cells.Item($r, 1).Value = "'" + $code
    $bowlingExtra.Cells.Item($r, 2).Value = "'" + $maidenOvers
    $bowlingExtra.Cells.Item($r, 3).Value = "'" + $percentWickets
}

# The leading apostrophes force everything to be stored as text (matching
# the source data, which keeps match codes / stats as strings); strip the
# resulting "quote prefix" styling so the cells stay on the default style,
# same as the rest of the workbook.
$bowlingExtra.Range("A2:C21").Style = "Normal"

Write-Host "Added 'ODI Bowling Extra' sheet and cleaned up 'ODI Batting Extra' row 2"
